# Auto-generated edit script applying cryptos list update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.262.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.74%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.693.98"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.28%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "520.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.96%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.71%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.993"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.62%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.571"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.10%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.730.21"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.56%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.28"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.49%  "

# Row 11
$ws.Range("E11").Value = "  +7.45%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.339"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.04%  "

# Row 13
$ws.Range("E13").Value = "  -0.83%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.164.88"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.18%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.196.23"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.64%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.02%  "

# Row 17
$ws.Range("E17").Value = "  +3.75%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.721.83"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.38%  "

# Row 19
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "353.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.51%  "

# Row 20
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.59"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.92%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.31%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.54%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.05%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.40%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.426"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.44%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.803.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.99%  "

# Row 27
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.162"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.92%  "

# Row 28
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.990"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.91%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0830"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.20%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.75%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.995"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.47%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.47"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +12.29%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.22"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.96%  "

# Row 34
$ws.Range("E34").Value = "  +3.76%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "150.32"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.86%  "

# Row 36
$ws.Range("E36").Value = "  +17.68%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.11"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.37%  "

# Row 38
$ws.Range("E38").Value = "  +5.70%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.867"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.89%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.97"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.17%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.74"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.08%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.43"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.65%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.629"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.00%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "284.13"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.68%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.77%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0988"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.78%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.991"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.72%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0538"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.01%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.80"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.62%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0233"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.56%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.019.03"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.08%  "
